$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.139097213745117
$ws.Range("C2").Value = 5.862069129943848
$ws.Range("D2").Value = 14.725564002990723
$ws.Range("E2").Value = 57.85714340209961
